# Apply the crypto price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / link / percentage cells can be written directly.
$ws.Range("D2").Value = '30.460.33'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '1.928.96'
$ws.Range("E3").Value = '  +4.10%  '
$ws.Range("E4").Value = '  -0.36%  '
$ws.Range("E5").Value = '  +3.06%  '
$ws.Range("E6").Value = '  -0.26%  '
$ws.Range("E7").Value = '  +0.50%  '
$ws.Range("E8").Value = '  +4.43%  '
$ws.Range("E9").Value = '  +4.34%  '
$ws.Range("E10").Value = '  +8.87%  '
$ws.Range("E11").Value = '  +25.95%  '
$ws.Range("D12").Value = '1.916.79'
$ws.Range("E12").Value = '  +3.47%  '
$ws.Range("E13").Value = '  +2.10%  '
$ws.Range("E14").Value = '  +3.41%  '
$ws.Range("E15").Value = '  +5.48%  '
$ws.Range("E16").Value = '  +23.73%  '
$ws.Range("D17").Value = '30.431.14'
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("E18").Value = '  +2.30%  '
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("E19").Value = '  -0.27%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("E20").Value = '  +2.64%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.167.62'
$ws.Range("E21").Value = '  +3.49%  '
$ws.Range("E22").Value = '  -0.24%  '
$ws.Range("E23").Value = '  +7.38%  '
$ws.Range("E24").Value = '  +6.84%  '
$ws.Range("E25").Value = '  +2.34%  '
$ws.Range("E26").Value = '  +1.17%  '
$ws.Range("E27").Value = '  +10.86%  '
$ws.Range("E28").Value = '  +7.25%  '
$ws.Range("E29").Value = '  +8.93%  '
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("E31").Value = '  +1.15%  '
$ws.Range("E32").Value = '  +2.65%  '
$ws.Range("E33").Value = '  +3.61%  '
$ws.Range("E34").Value = '  +6.65%  '
$ws.Range("E35").Value = '  +2.03%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("E36").Value = '  +1.10%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("E37").Value = '  +3.05%  '
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("E38").Value = '  +0.70%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("E39").Value = '  +2.91%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("B41").Value = 'Quant'
$ws.Range("C41").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("E41").Value = '  +0.75%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("E42").Value = '  +4.99%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("E43").Value = '  +10.72%  '
$ws.Range("E44").Value = '  -0.36%  '
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("E45").Value = '  +2.33%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("E46").Value = '  +0.59%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("E47").Value = '  +7.63%  '
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("E48").Value = '  +3.40%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("E49").Value = '  +0.83%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("E50").Value = '  +1.92%  '
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("E51").Value = '  +4.39%  '

# Numeric-looking price strings (e.g. "1.000", "5.800") must stay as literal
# text, matching the original inline strings, instead of being coerced to
# numbers by Excel. Stage them in a scratch cell formatted as Text, then
# copy/paste-special (values only) into the target cell so no number
# conversion / stray cell formatting is applied to the target cell.
$staging = $ws.Range("Z1")
$staging.NumberFormat = "@"
$staging.Value = '0.9978'
$staging.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$staging.Value = '240.74'
$staging.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$staging.Value = '0.9988'
$staging.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$staging.Value = '0.4769'
$staging.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$staging.Value = '0.2872'
$staging.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$staging.Value = '0.06598'
$staging.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$staging.Value = '19.17'
$staging.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$staging.Value = '106.61'
$staging.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$staging.Value = '5.133'
$staging.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$staging.Value = '0.6592'
$staging.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$staging.Value = '305.06'
$staging.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$staging.Value = '12.98'
$staging.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$staging.Value = '1.000'
$staging.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$staging.Value = '0.000007512'
$staging.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$staging.Value = '1.000'
$staging.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$staging.Value = '5.280'
$staging.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$staging.Value = '6.315'
$staging.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$staging.Value = '168.64'
$staging.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$staging.Value = '9.230'
$staging.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$staging.Value = '19.93'
$staging.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$staging.Value = '2.007'
$staging.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$staging.Value = '0.1122'
$staging.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$staging.Value = '4.094'
$staging.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$staging.Value = '3.932'
$staging.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$staging.Value = '0.05014'
$staging.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$staging.Value = '0.7425'
$staging.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$staging.Value = '1.152'
$staging.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$staging.Value = '2.728'
$staging.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$staging.Value = '0.01951'
$staging.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$staging.Value = '2.699'
$staging.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$staging.Value = '2.059'
$staging.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$staging.Value = '0.8751'
$staging.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$staging.Value = '107.19'
$staging.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$staging.Value = '5.800'
$staging.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$staging.Value = '69.78'
$staging.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$staging.Value = '0.9987'
$staging.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$staging.Value = '0.4144'
$staging.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$staging.Value = '7.229'
$staging.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$staging.Value = '9.220'
$staging.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$staging.Value = '34.83'
$staging.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$staging.Value = '0.1206'
$staging.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$staging.Value = '0.05626'
$staging.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$staging.Value = '0.3841'
$staging.Copy()
$ws.Range("D51").PasteSpecial(-4163)

$staging.Clear()
$excel.CutCopyMode = 0
